$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set: years 2010-2022 with corresponding birth/death/natural-growth rates
$data = @(
    @("2010年", 11.9,  7.11, 4.79),
    @("2011年", 13.27, 7.14, 6.13),
    @("2012年", 14.57, 7.13, 7.43),
    @("2013年", 13.03, 7.13, 5.9),
    @("2014年", 13.83, 7.12, 6.71),
    @("2015年", 11.99, 7.07, 4.93),
    @("2016年", 13.57, 7.04, 6.53),
    @("2017年", 12.64, 7.06, 5.58),
    @("2018年", 10.86, 7.08, 3.78),
    @("2019年", 10.41, 7.09, 3.32),
    @("2020年", 8.52,  7.07, 1.45),
    @("2021年", 7.52,  7.18, 0.34),
    @("2022年", 6.77,  7.37, -0.6)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# Remove the now-unused trailing rows (previously through row 22)
$lastOldRow = 22
$firstExtraRow = $data.Count + 2
$ws.Range("A" + $firstExtraRow + ":D" + $lastOldRow).EntireRow.Delete()
